$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1358.4814
$ws.Range("I15").Value = 1358.4814
$ws.Range("K15").Value = 4075.4442
$ws.Range("M15").Value = -3906.4442
$ws.Range("H28").Value = 10008.546
$ws.Range("I28").Value = 12838
$ws.Range("K28").Value = 12838
$ws.Range("M28").Value = -12353
$ws.Range("H76").Value = 3925.889
$ws.Range("I76").Value = 3458.25
$ws.Range("J76").Value = 4300
$ws.Range("K76").Value = 3458.25
$ws.Range("L76").Value = 4300
$ws.Range("M76").Value = -3143.25
$ws.Range("N76").Value = -4930
$ws.Range("H79").Value = 3925.889
$ws.Range("I79").Value = 3458.25
$ws.Range("J79").Value = 4300
$ws.Range("K79").Value = 3458.25
$ws.Range("L79").Value = 4300
$ws.Range("M79").Value = -2366.25
$ws.Range("N79").Value = -6484
$ws.Range("H87").Value = 44999
$ws.Range("J87").Value = 44999
$ws.Range("L87").Value = 44999
$ws.Range("N87").Value = -47495
$ws.Range("H90").Value = 44999
$ws.Range("J90").Value = 44999
$ws.Range("L90").Value = 134997
$ws.Range("N90").Value = -147477
$ws.Range("H106").Value = 4914
$ws.Range("I106").Value = 4914
$ws.Range("K106").Value = 4914
$ws.Range("M106").Value = -4283
$ws.Range("H112").Value = 3191.3333
$ws.Range("J112").Value = 4124.625
$ws.Range("L112").Value = 12373.875
$ws.Range("N112").Value = -14589.875
$ws.Range("H137").Value = 2992.8235
$ws.Range("I137").Value = 1428.4286
$ws.Range("K137").Value = 4285.2858
$ws.Range("M137").Value = -1735.2858

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6659.85
$ws.Range("I32").Value = 6659.85
$ws.Range("K32").Value = 6659.85
$ws.Range("M32").Value = -6372.85
$ws.Range("H61").Value = 3581.6667
$ws.Range("I61").Value = 3498
$ws.Range("K61").Value = 3498
$ws.Range("M61").Value = -3286
$ws.Range("H74").Value = 4770.263
$ws.Range("J74").Value = 6199.6
$ws.Range("L74").Value = 6199.6
$ws.Range("N74").Value = -7947.6
$ws.Range("H77").Value = 4770.263
$ws.Range("J77").Value = 6199.6
$ws.Range("L77").Value = 30998
$ws.Range("N77").Value = -39734
$ws.Range("H110").Value = 2200
$ws.Range("I110").Value = 2288.5715
$ws.Range("K110").Value = 2288.5715
$ws.Range("M110").Value = -243.5715
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1548.25
$ws.Range("I132").Value = 1568.7407
$ws.Range("K132").Value = 4706.2221
$ws.Range("M132").Value = -2176.2221
$ws.Range("H136").Value = 3581.6667
$ws.Range("I136").Value = 3498
$ws.Range("K136").Value = 10494
$ws.Range("M136").Value = -7944

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4892
$ws.Range("I20").Value = 4408.769
$ws.Range("J20").Value = 5939
$ws.Range("K20").Value = 4408.769
$ws.Range("L20").Value = 5939
$ws.Range("M20").Value = -4161.769
$ws.Range("N20").Value = -6433
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227
$ws.Range("H99").Value = 2774.25
$ws.Range("I99").Value = 3028.4285
$ws.Range("K99").Value = 3028.4285
$ws.Range("M99").Value = -1530.4285
$ws.Range("H106").Value = 17906.166
$ws.Range("J106").Value = 17906.166
$ws.Range("L106").Value = 17906.166
$ws.Range("N106").Value = -20430.166
$ws.Range("H134").Value = 3449.56
$ws.Range("I134").Value = 3493.3333
$ws.Range("K134").Value = 10479.9999
$ws.Range("M134").Value = -7944.999899999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2446.279
$ws.Range("I31").Value = 1992.129
$ws.Range("K31").Value = 1992.129
$ws.Range("M31").Value = -1697.129
$ws.Range("H34").Value = 2446.279
$ws.Range("I34").Value = 1992.129
$ws.Range("K34").Value = 1992.129
$ws.Range("M34").Value = -1790.129
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920
$ws.Range("H132").Value = 999.7
$ws.Range("I132").Value = 999.7
$ws.Range("K132").Value = 2999.1
$ws.Range("M132").Value = -469.1000000000004
$ws.Range("H134").Value = 2323.1
$ws.Range("I134").Value = 2359
$ws.Range("K134").Value = 7077
$ws.Range("M134").Value = -4542

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1654.4286
$ws.Range("I5").Value = 1231
$ws.Range("J5").Value = 1972
$ws.Range("K5").Value = 3693
$ws.Range("L5").Value = 5916
$ws.Range("M5").Value = -3581
$ws.Range("N5").Value = -6140
$ws.Range("H129").Value = 543.75
$ws.Range("I129").Value = 514.2857
$ws.Range("K129").Value = 1542.8571
$ws.Range("M129").Value = 3457.1429
$ws.Range("H131").Value = 1200.7142
$ws.Range("J131").Value = 1581.8
$ws.Range("L131").Value = 4745.4
$ws.Range("N131").Value = -14825.4
$ws.Range("H135").Value = 1654.4286
$ws.Range("I135").Value = 1231
$ws.Range("J135").Value = 1972
$ws.Range("K135").Value = 11079
$ws.Range("L135").Value = 17748
$ws.Range("M135").Value = -8544
$ws.Range("N135").Value = -22818

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2873.5715
$ws.Range("I80").Value = 1625
$ws.Range("J80").Value = 3810
$ws.Range("K80").Value = 1625
$ws.Range("L80").Value = 3810
$ws.Range("M80").Value = -627
$ws.Range("N80").Value = -5806
$ws.Range("H83").Value = 2873.5715
$ws.Range("I83").Value = 1625
$ws.Range("J83").Value = 3810
$ws.Range("K83").Value = 8125
$ws.Range("L83").Value = 19050
$ws.Range("M83").Value = -3133
$ws.Range("N83").Value = -29034
$ws.Range("H107").Value = 975
$ws.Range("J107").Value = 3000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840
$ws.Range("H132").Value = 4337.3335
$ws.Range("I132").Value = 4006
$ws.Range("K132").Value = 12018
$ws.Range("M132").Value = -9488

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1652.1818
$ws.Range("I82").Value = 972.25
$ws.Range("K82").Value = 972.25
$ws.Range("M82").Value = -611.25
$ws.Range("H85").Value = 1652.1818
$ws.Range("I85").Value = 972.25
$ws.Range("K85").Value = 972.25
$ws.Range("M85").Value = 275.75
$ws.Range("H132").Value = 4890.5264
$ws.Range("I132").Value = 3856.111
$ws.Range("J132").Value = 5821.5
$ws.Range("K132").Value = 11568.333
$ws.Range("L132").Value = 17464.5
$ws.Range("M132").Value = -9038.332999999999
$ws.Range("N132").Value = -22524.5
$ws.Range("H136").Value = 3502.3333
$ws.Range("J136").Value = 3500
$ws.Range("L136").Value = 10500
$ws.Range("N136").Value = -15600

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 507500
$ws.Range("J29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15580
$ws.Range("H80").Value = 20301
$ws.Range("J80").Value = 20301
$ws.Range("L80").Value = 20301
$ws.Range("N80").Value = -22297
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766
$ws.Range("H83").Value = 20301
$ws.Range("J83").Value = 20301
$ws.Range("L83").Value = 60903
$ws.Range("N83").Value = -70887
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652
$ws.Range("H122").Value = 4502.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4502.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13507.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -18407.5
$ws.Range("H131").Value = 100715
$ws.Range("J131").Value = 100715
$ws.Range("L131").Value = 100715
$ws.Range("N131").Value = -110795
$ws.Range("H136").Value = 9866.444
$ws.Range("J136").Value = 10400
$ws.Range("L136").Value = 31200
$ws.Range("N136").Value = -36300
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280

Write-Host "Applied all Marilith_Profits updates"